$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Tests": insert a new "param:scope" column before the existing
# "param:q" column (I), and add a new row 3 for the "Missing Required Param"
# test case.
# ---------------------------------------------------------------------------
$tests = $wb.Worksheets.Item("Tests")

# Insert a new column at I, shifting param:q..param:synonymsEnabled right by one.
$tests.Columns("I:I").Insert()

# Header for the new column.
$tests.Range("I1").Value = "param:scope"

# New value for existing row 2 (Basic Test).
$tests.Range("I2").Value = "work"

# New row 3: "Missing Required Param" test case.
# Note: values that look like booleans/numbers are prefixed with a leading
# apostrophe so Excel stores them as literal text (matching the workbook's
# existing text-typed "true"/"1"/"20" cells) instead of auto-coercing them
# to boolean/numeric types. The style is reset to "Normal" afterwards so the
# quote-prefix formatting flag doesn't linger on the cell.
$tests.Range("A3").Value = "get-search-estimate - Missing Required Param"
$tests.Range("B3").Value = "Test GET /api/search-estimate/:scope with missing required parameters"
$tests.Range("C3").Value = "'true"
$tests.Range("C3").Style = "Normal"
$tests.Range("D3").Value = 400
$tests.Range("E3").Value = 10000
$tests.Range("F3").Value = 2000
$tests.Range("G3").Value = 500
$tests.Range("H3").Value = "get-search-estimate,validation"
$tests.Range("I3").Value = "'"
$tests.Range("I3").Style = "Normal"
$tests.Range("J3").Value = "test query"
$tests.Range("K3").Value = "'true"
$tests.Range("K3").Style = "Normal"
$tests.Range("L3").Value = "'1"
$tests.Range("L3").Style = "Normal"
$tests.Range("M3").Value = "'20"
$tests.Range("M3").Style = "Normal"
$tests.Range("N3").Value = "example"
$tests.Range("O3").Value = "example"
$tests.Range("P3").Value = "'true"
$tests.Range("P3").Style = "Normal"
$tests.Range("Q3").Value = "'true"
$tests.Range("Q3").Style = "Normal"
$tests.Range("R3").Value = "'true"
$tests.Range("R3").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "Documentation": insert a "param:scope" description row before the
# existing "param:q" row, and a "Required parameters" note row before the
# "Optional parameters" note row.
# ---------------------------------------------------------------------------
$docs = $wb.Worksheets.Item("Documentation")

# Insert a row at 18 for the new "param:scope" description.
$docs.Rows("18:18").Insert()
$docs.Range("A18").Value = "param:scope"
$docs.Range("B18").Value = "Search scope (work, person, place, concept, event, etc.) (string) (REQUIRED - highlighted in yellow)"

# Insert a row at 33 (after the insert above, row 33 holds the "Optional
# parameters" note) for the new "Required parameters" note.
$docs.Rows("33:33").Insert()
$docs.Range("A33").Value = "• Required parameters: scope"

Write-Output "done"
